# Edit slide 12 ("Probe Query for Links") of the TWAMP SRPM deck.
#
# 1. Resize/reposition the "Content Placeholder 2" text box (shrink its
#    width and nudge it right/down), and drop the second paragraph
#    ("For DM, payload contains RFC 5357 ...") while turning off the
#    bullet on what becomes the trailing empty paragraph.
# 2. Nudge the "Rectangle 4" diagram box up slightly (Top only).
# 3. Blank out the " for Delay Measurement" suffix in the diagram text,
#    replacing it with spaces so the line length is preserved.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)

# --- Shape: "Content Placeholder 2" --------------------------------
$contentShape = $s.Shapes.Item(3)

# Reposition / resize (EMU -> points, 12700 EMU per point).
$contentShape.Left = 48
$contentShape.Top = 65.2115
$contentShape.Width = 624
$contentShape.Height = 67.5

$tr = $contentShape.TextFrame.TextRange
# Remove the second paragraph ("For DM, payload contains RFC 5357 ...");
# its text merges away, leaving the trailing empty paragraph in place.
$para2 = $tr.Paragraphs(2, 1)
$para2.Delete()

# Turn off the bullet on the (now second / last) empty paragraph.
$tr = $contentShape.TextFrame.TextRange
$lastPara = $tr.Paragraphs(2, 1)
$lastPara.ParagraphFormat.Bullet.Visible = 0

# --- Shape: "Rectangle 4" -------------------------------------------
$rectShape = $s.Shapes.Item(4)
$rectShape.Top = 127.8939

# Blank out " for Delay Measurement" (replace with equal-length spaces).
$rectTr = $rectShape.TextFrame.TextRange
$found = $rectTr.Find("User-configured Port for Delay Measurement", 0)
$found.Text = "User-configured Port                      "
